# Vega Modelo de Temuco - Espárragos: weekly refresh.
# Two new daily observations are inserted into the historical series:
#   - a new row at position 12 (pushes the old row 12.. down by one)
#   - a new row at position 26 (pushes the old row 26.. down by one more)
# The net effect (matching the authored diff) is that the sheet grows from
# A1:R73 to A1:R75, with every existing record below row 11 shifted down by
# one or two rows, and the two brand new records written into the freshly
# inserted rows 12 and 26.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shared / constant columns for every data row in this sheet.
$mercadoId = 10
$mercado   = "Vega Modelo de Temuco"
$region    = "La Araucanía"
$codreg    = 9
$catId     = 300000000
$categoria = "Espárragos"
$clasif    = "Hortaliza"

# --- Insert the first new row at 12 ------------------------------------
$ws.Rows(12).Insert()

$ws.Range("A12").Value = $mercadoId
$ws.Range("B12").Value = $mercado
$ws.Range("C12").Value = $region
$ws.Range("D12").Value = 44847
$ws.Range("E12").Value = $codreg
$ws.Range("F12").Value = $catId
$ws.Range("G12").Value = $categoria
$ws.Range("H12").Value = "Sin especificar"
$ws.Range("I12").Value = "Primera"
$ws.Range("J12").Value = 900
$ws.Range("K12").Value = 1500
$ws.Range("L12").Value = 1600
$ws.Range("M12").Value = 1533
$ws.Range("N12").Value = '$/kilo'
$ws.Range("O12").Value = "Región del Maule"
$ws.Range("P12").Value = 1533
$ws.Range("Q12").Value = 1
$ws.Range("R12").Value = $clasif

# --- Insert the second new row at 26 ------------------------------------
$ws.Rows(26).Insert()

$ws.Range("A26").Value = $mercadoId
$ws.Range("B26").Value = $mercado
$ws.Range("C26").Value = $region
$ws.Range("D26").Value = 44848
$ws.Range("E26").Value = $codreg
$ws.Range("F26").Value = $catId
$ws.Range("G26").Value = $categoria
$ws.Range("H26").Value = "Sin especificar"
$ws.Range("I26").Value = "Primera"
$ws.Range("J26").Value = 550
$ws.Range("K26").Value = 1500
$ws.Range("L26").Value = 1600
$ws.Range("M26").Value = 1536
$ws.Range("N26").Value = '$/kilo'
$ws.Range("O26").Value = "Región del Maule"
$ws.Range("P26").Value = 1536
$ws.Range("Q26").Value = 1
$ws.Range("R26").Value = $clasif

Write-Host "Done. Dimension now $($ws.UsedRange.Rows.Count) rows."
